# Make address optional for paatokset
# Clear the address-related cells (N3/O3/P3: street, postal code, city) on row 3
# so the "Kohteen osoite" address columns are empty/optional, matching the
# already-blank treatment used elsewhere (e.g. row 6).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N3").Value = ""
$ws.Range("O3").Value = ""
$ws.Range("P3").Value = ""

# Move the active selection to A4 (matches the recorded selection state).
$ws.Range("A4").Select()
